# Added Enquiry to remaining Params
# Adds a "Menu Option" legend/help block in columns Z:AA (rows 2-5) describing
# the menu option fields, formatted with a thin box border, wrapped + vertically
# centered text, and a slightly smaller (7.5pt) font.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column widths for the added legend columns ---
$ws.Columns.Item(26).ColumnWidth = 51.92   # column Z  -> ~52.78 chars
$ws.Columns.Item(27).ColumnWidth = 46.75   # column AA -> ~47.66 chars

# --- New legend content ---
$ws.Range("Z2").Value  = "Menu Option"
$ws.Range("AA2").Value = "Menu Option and related operation"
$ws.Range("Z3").Value  = "Brief Description of the Option"
$ws.Range("Z4").Value  = "URL of Menu Option"
$ws.Range("Z5").Value  = "Transaction Code"

# --- Formatting: thin box border, vertical-centered + wrapped 7.5pt text ---
# Build the format once on Z2, then fan it out to the other touched cells with
# a format-only paste so they all share the same single style record.
$base = $ws.Range("Z2")
$style = $base.Style
$style.Font.Size = 7.5
$style.Borders.LineStyle = 1
$style.VerticalAlignment = -4108
$style.WrapText = $true

$base.Copy()
$ws.Range("AA2").PasteSpecial(-4122)
$ws.Range("Z3:Z5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection / view state to match the authored workbook ---
$ws.Range("Z2:AA5").Select()
$excel.ActiveWindow.ScrollColumn = 16
